$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: Israel's national mitigation target row loses its Topic/Unit/Scale/Time/Principle/Notes
# columns and its "Relevance" flips from yes -> no.
$ws.Range("B12").Value = "no"
$ws.Range("C12:H12").ClearContents()

# Row 37: the "other(knowledge)" / redistribution row is dropped entirely (Relevance -> no)
$ws.Range("B37").Value = "no"
$ws.Range("C37:H37").ClearContents()

# Row 38: the "other(innovations), other(technology)" / sharing row is dropped entirely (Relevance -> no)
$ws.Range("B38").Value = "no"
$ws.Range("C38:H38").ClearContents()

# Row 54: "other(cooperation)" -> "cooperation", "utilitarian" -> "utilitarian, egalitarian"
$ws.Range("G54").Value = "utilitarian, egalitarian"
$ws.Range("C54").Value = "cooperation"

# Row 27: "other(support)" -> "support"
$ws.Range("D27").Value = "support"

# Update the active selection/view to rest on D27 (matches the author's saved view state)
[void]$ws.Range("D27").Select()
